$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A7").EntireRow.Insert()
